$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) - Worksheets.Item(1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1524
$ws.Range("F3").Value = 862
$ws.Range("F4").Value = 446
$ws.Range("F5").Value = 897
$ws.Range("F6").Value = 507
$ws.Range("F7").Value = 7622
$ws.Range("F10").Value = 1936
$ws.Range("F11").Value = 5520
$ws.Range("F12").Value = 567
$ws.Range("F14").Value = 7611
$ws.Range("F15").Value = 8987
$ws.Range("F16").Value = 1145
$ws.Range("F17").Value = 897
$ws.Range("F18").Value = 4448
$ws.Range("F19").Value = 671
$ws.Range("F20").Value = 228
$ws.Range("F22").Value = 283
$ws.Range("F24").Value = 1194
$ws.Range("F25").Value = 112
$ws.Range("F26").Value = 1658
$ws.Range("F27").Value = 716
$ws.Range("F28").Value = 925
$ws.Range("F30").Value = 1872
$ws.Range("F31").Value = 335
$ws.Range("F32").Value = 2290
$ws.Range("F34").Value = 111
$ws.Range("F35").Value = 1461
$ws.Range("F38").Value = 793
$ws.Range("F40").Value = 2971
$ws.Range("F41").Value = 4098
$ws.Range("F43").Value = 39
$ws.Range("F44").Value = 416
$ws.Range("F45").Value = 509
$ws.Range("F48").Value = 171
$ws.Range("F49").Value = 4084

# 演出 (Performance) - Worksheets.Item(2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F27").Value = 35

# 本地生活 (Local life) - Worksheets.Item(3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5224

# 全部类型 (All types) - Worksheets.Item(4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1524
$ws.Range("F4").Value = 862
$ws.Range("F5").Value = 446
$ws.Range("F6").Value = 897
$ws.Range("F7").Value = 507
$ws.Range("F11").Value = 5520
$ws.Range("F12").Value = 567
$ws.Range("F13").Value = 7611
$ws.Range("F15").Value = 1145
$ws.Range("F16").Value = 897
$ws.Range("F17").Value = 671
$ws.Range("F18").Value = 228
$ws.Range("F20").Value = 283
$ws.Range("F23").Value = 1194
$ws.Range("F24").Value = 112
$ws.Range("F25").Value = 1658
$ws.Range("F26").Value = 716
$ws.Range("F27").Value = 925
$ws.Range("F29").Value = 1872
$ws.Range("F30").Value = 335
$ws.Range("F31").Value = 2290
$ws.Range("F40").Value = 4098
$ws.Range("F41").Value = 35
$ws.Range("F43").Value = 39
$ws.Range("F44").Value = 416
$ws.Range("F45").Value = 509
$ws.Range("F48").Value = 171
$ws.Range("F49").Value = 4084
